$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "56.046.14"
$ws.Range("E2").Value = "  -2.67%  "
$ws.Range("D3").Value = "2.367.48"
$ws.Range("E3").Value = "  -1.56%  "
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "503.75"
$ws.Range("E5").Value = "  -0.33%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "131.16"
$ws.Range("E6").Value = "  -1.20%  "
$ws.Range("E7").Value = "  +0.65%  "
$ws.Range("E8").Value = "  -2.04%  "
$ws.Range("D9").Value = "2.368.83"
$ws.Range("E9").Value = "  -3.00%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0973"
$ws.Range("E10").Value = "  -0.31%  "
$ws.Range("E11").Value = "  +0.51%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.328"
$ws.Range("E12").Value = "  +1.34%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.68"
$ws.Range("E13").Value = "  +0.31%  "
$ws.Range("D14").Value = "2.786.46"
$ws.Range("E14").Value = "  -1.43%  "
$ws.Range("D15").Value = "55.980.61"
$ws.Range("E15").Value = "  -2.54%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "21.45"
$ws.Range("E16").Value = "  -2.14%  "
$ws.Range("E17").Value = "  -1.55%  "
$ws.Range("D18").Value = "2.406.73"
$ws.Range("E18").Value = "  -0.85%  "
$ws.Range("E19").Value = "  -2.73%  "
$ws.Range("E20").Value = "  -2.96%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "306.97"
$ws.Range("E21").Value = "  -2.38%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.25"
$ws.Range("E22").Value = "  -1.94%  "
$ws.Range("E23").Value = "  +0.40%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "65.08"
$ws.Range("E24").Value = "  -0.88%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.00"
$ws.Range("E25").Value = "  +1.04%  "
$ws.Range("E26").Value = "  -3.79%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.148"
$ws.Range("E27").Value = "  -3.03%  "
$ws.Range("E28").Value = "  -3.80%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "172.77"
$ws.Range("E29").Value = "  -0.22%  "
$ws.Range("E30").Value = "  -2.68%  "
$ws.Range("E31").Value = "  -2.71%  "
$ws.Range("E32").Value = "  +0.20%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.10"
$ws.Range("E33").Value = "  -4.65%  "
$ws.Range("E34").Value = "  -7.47%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.997"
$ws.Range("E35").Value = "  +0.56%  "
$ws.Range("E36").Value = "  -3.22%  "
$ws.Range("E37").Value = "  -3.83%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.79"
$ws.Range("E38").Value = "  -1.69%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "36.02"
$ws.Range("E39").Value = "  -0.57%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.797"
$ws.Range("E40").Value = "  -2.58%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.42"
$ws.Range("E41").Value = "  -3.15%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "131.29"
$ws.Range("E42").Value = "  -1.55%  "
$ws.Range("E43").Value = "  -1.82%  "
$ws.Range("E44").Value = "  -5.36%  "
$ws.Range("E45").Value = "  -1.16%  "
$ws.Range("E46").Value = "  -0.82%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "244.23"
$ws.Range("E47").Value = "  -6.38%  "
$ws.Range("E48").Value = "  -3.42%  "
$ws.Range("E49").Value = "  -2.75%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "17.04"
$ws.Range("E50").Value = "  -1.19%  "
$ws.Range("B51").Value = "dogwifhat"
$ws.Range("C51").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.55"
$ws.Range("E51").Value = "  -2.29%  "
